$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Update selection
$ws.Range("E16").Select()

# Update passwords for rows 2-7 (Business/Manufacturer/AuthorisedRep for Auto & Noor.Uddin)
$ws.Range("B2").Value = "MHRA12345"
$ws.Range("B3").Value = "MHRA12345"
$ws.Range("B4").Value = "MHRA12345"
$ws.Range("B5").Value = "MHRA12345"
$ws.Range("B6").Value = "MHRA12345"
$ws.Range("B7").Value = "MHRA12345"

# Update ignore flag for rows 14-16 (Lambros.Poullais Business/Manufacturer/AuthorisedRep)
$ws.Range("C14").Value = "no"
$ws.Range("C15").Value = "no"
$ws.Range("C16").Value = "no"
